# Fixed issue with Analytics not working
# ------------------------------------------------------------------
# This reproduces the genuine semantic edits behind the diff:
#   1. The date/time numeric format used for column A header rows
#      (style index 2, numFmt 164) is changed from
#      "YYYY-MM-DD HH:MM:SS" to "yyyy\-mm\-dd\ hh:mm:ss".
#   2. Column A on both worksheets is given an explicit custom width
#      (a user resized the column).
#   3. The view state on ADAM_MBM_Worked is scrolled down and a new
#      range (A56:G58) is selected, reflecting where the user was
#      working when the sheet was last saved.
# (All other differences in the raw OOXML -- refreshed shared-string
#  ordering, namespace/schema upgrades, fileVersion/calcId bumps,
#  theme font-list additions, etc. -- are byproducts of the workbook
#  being re-saved by a newer Excel build and carry no data meaning;
#  they are not something a user action produces through the object
#  model.)
# ------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ADAM_MBM_Worked")
$ws2 = $wb.Worksheets.Item("ADAM_UET_Worked")

# 1. Update the date/time display format applied to the timestamp
#    cells in column A (both sheets use the same style).
$ws1.Range("A2:A3").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws2.Range("A2").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# 2. Resize column A on each sheet.
$ws1.Columns.Item(1).ColumnWidth = 33.71
$ws2.Columns.Item(1).ColumnWidth = 19.71

# 3. Restore the on-screen scroll position / selection that was
#    active on ADAM_MBM_Worked.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 49
$ws1.Range("A56:G58").Select()
